# Updates the cryptos list (coin name/link/price/volume) to the latest
# scraped snapshot. A handful of rows (20/21, 42/43, 47/48) also swapped
# rank order, so their Coin/Link columns are rewritten along with
# Price/Volume(1h).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; B=$null; C=$null; D='27.503.02'; E='  +2.10%  ' },
    @{ Row=3; B=$null; C=$null; D='1.866.95'; E='  +1.15%  ' },
    @{ Row=4; B=$null; C=$null; D='1.014'; E='  +0.30%  ' },
    @{ Row=5; B=$null; C=$null; D='311.75'; E='  +0.67%  ' },
    @{ Row=7; B=$null; C=$null; D='0.4777'; E='  -0.12%  ' },
    @{ Row=8; B=$null; C=$null; D='0.3728'; E='  +1.62%  ' },
    @{ Row=9; B=$null; C=$null; D='0.07286'; E='  +0.84%  ' },
    @{ Row=10; B=$null; C=$null; D='0.9341'; E='  +0.82%  ' },
    @{ Row=11; B=$null; C=$null; D='20.67'; E='  +5.03%  ' },
    @{ Row=12; B=$null; C=$null; D='0.07831'; E='  +1.59%  ' },
    @{ Row=13; B=$null; C=$null; D='1.873.31'; E='  +2.06%  ' },
    @{ Row=14; B=$null; C=$null; D='5.437'; E='  +2.26%  ' },
    @{ Row=15; B=$null; C=$null; D='6.544'; E='  +2.11%  ' },
    @{ Row=16; B=$null; C=$null; D='90.14'; E='  +1.49%  ' },
    @{ Row=17; B=$null; C=$null; D=$null; E='  +0.24%  ' },
    @{ Row=18; B=$null; C=$null; D='0.000008888'; E='  +2.83%  ' },
    @{ Row=19; B=$null; C=$null; D=$null; E='  +0.19%  ' },
    @{ Row=20; B='WrappedBTC'; C='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D='27.523.94'; E='  +2.04%  ' },
    @{ Row=21; B='Avalanche'; C='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D='14.51'; E='  -0.27%  ' },
    @{ Row=22; B=$null; C=$null; D='5.114'; E='  +1.07%  ' },
    @{ Row=23; B=$null; C=$null; D='10.70'; E='  +0.45%  ' },
    @{ Row=24; B=$null; C=$null; D='1.949'; E='  +1.12%  ' },
    @{ Row=26; B=$null; C=$null; D='18.44'; E='  +1.50%  ' },
    @{ Row=27; B=$null; C=$null; D='2.018'; E='  +1.22%  ' },
    @{ Row=28; B=$null; C=$null; D='115.64'; E='  +1.31%  ' },
    @{ Row=29; B=$null; C=$null; D='4.989'; E='  +0.94%  ' },
    @{ Row=30; B=$null; C=$null; D='0.08910'; E='  +0.29%  ' },
    @{ Row=31; B=$null; C=$null; D='3.335'; E='  +0.54%  ' },
    @{ Row=32; B=$null; C=$null; D='1.216'; E='  +3.74%  ' },
    @{ Row=33; B=$null; C=$null; D='0.7569'; E='  +1.73%  ' },
    @{ Row=34; B=$null; C=$null; D='4.611'; E='  +2.65%  ' },
    @{ Row=35; B=$null; C=$null; D='2.721'; E='  +0.06%  ' },
    @{ Row=36; B=$null; C=$null; D='0.02048'; E='  +4.74%  ' },
    @{ Row=37; B=$null; C=$null; D='1.120'; E='  -0.70%  ' },
    @{ Row=38; B=$null; C=$null; D='3.004'; E='  +0.69%  ' },
    @{ Row=39; B=$null; C=$null; D='0.05268'; E='  +0.07%  ' },
    @{ Row=40; B=$null; C=$null; D='0.5319'; E='  +2.31%  ' },
    @{ Row=41; B=$null; C=$null; D='7.061'; E='  +0.95%  ' },
    @{ Row=42; B='Algorand'; C='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; D='0.1523'; E='  +0.84%  ' },
    @{ Row=43; B='Aptos'; C='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D='8.474'; E='  +3.34%  ' },
    @{ Row=44; B=$null; C=$null; D='10.71'; E='  +1.73%  ' },
    @{ Row=45; B=$null; C=$null; D='0.4792'; E='  +1.38%  ' },
    @{ Row=46; B=$null; C=$null; D=$null; E='  +0.24%  ' },
    @{ Row=47; B='Quant'; C='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D='102.75'; E='  +1.32%  ' },
    @{ Row=48; B='NEARProtocol'; C='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D='1.652'; E='  +3.11%  ' },
    @{ Row=49; B=$null; C=$null; D='67.31'; E='  +2.72%  ' },
    @{ Row=50; B=$null; C=$null; D='0.06095'; E='  +1.11%  ' },
    @{ Row=51; B=$null; C=$null; D='0.9180'; E='  +3.63%  ' }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($u.B -ne $null) {
        $ws.Cells.Item($row, 2).Value = $u.B
    }
    if ($u.C -ne $null) {
        $ws.Cells.Item($row, 3).Value = $u.C
    }
    if ($u.D -ne $null) {
        $priceCell = $ws.Cells.Item($row, 4)
        # Price column is stored as text (e.g. "27.503.02", "1.014"); many
        # values would otherwise be auto-coerced to a number by Excel, so
        # force a text format for the write, then drop back to the default
        # "Normal" style so no stray formatting is left behind.
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $u.D
        $priceCell.Style = "Normal"
    }
    if ($u.E -ne $null) {
        $ws.Cells.Item($row, 5).Value = $u.E
    }
}
